$d = $word.ActiveDocument

# Insert a brand-new empty paragraph right after the paragraph that ends
# with "...transports a DNA payload." (and before the following blank
# paragraph).
$anchor = $d.Content
$anchor.Find.Execute("transports a DNA payload.")
$anchor.Collapse(0)            # wdCollapseEnd
$anchor.InsertParagraphAfter()

# Locate the paragraph that still ends with "...payload." so we can find
# the freshly-inserted (empty) paragraph that immediately follows it.
$idx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*transports a DNA payload.*") {
        $idx = $i
        break
    }
}

$newPara = $d.Paragraphs.Item($idx + 1)
$r = $newPara.Range

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">In vivo gene editing: to correct mutations that causes genetic disease (or any other ethical genetic mutation): site-specific-cleavage was realized using </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:t>inc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> Finger Nucleases (ZFNs) and Transcription-Activator-like Effector Nucleases (TALENs). These techniques have been since then replaced by the more precise CRISPR-Cas enzyme which can cut DNA within a long target sequence (~30 bp) and can be more easily be reprogrammed for new targets. </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:t>The CRISPR-Cas</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:t>9</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> uses a small piece of RNA to recognize the intended DNA sequence, and guide the Cas 9 enzyme to cut the DNA at the targeted sequence.</w:t></w:r></w:p>'

$r.InsertXML($xml)

Write-Output "Inserted new paragraph after index $idx"
